$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.583.70"
$ws.Range("E2").Value = "  -7.07%  "
$ws.Range("D3").Value = "1.692.07"
$ws.Range("E3").Value = "  -5.54%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'219.92"
$ws.Range("E5").Value = "  -5.01%  "
$ws.Range("D6").Value = "'0.5122"
$ws.Range("E6").Value = "  -12.79%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "'0.2660"
$ws.Range("E8").Value = "  -3.68%  "
$ws.Range("D9").Value = "'22.05"
$ws.Range("E9").Value = "  -4.67%  "
$ws.Range("D10").Value = "'0.06312"
$ws.Range("E10").Value = "  -6.24%  "
$ws.Range("E11").Value = "  -2.02%  "
$ws.Range("D12").Value = "1.695.50"
$ws.Range("E12").Value = "  -5.40%  "
$ws.Range("D13").Value = "'4.518"
$ws.Range("E13").Value = "  -5.66%  "
$ws.Range("D14").Value = "'0.5785"
$ws.Range("E14").Value = "  -5.63%  "
$ws.Range("D15").Value = "1.923.12"
$ws.Range("E15").Value = "  -5.48%  "
$ws.Range("D16").Value = "'0.000008525"
$ws.Range("E16").Value = "  -4.17%  "
$ws.Range("D17").Value = "'65.31"
$ws.Range("E17").Value = "  -13.22%  "
$ws.Range("D18").Value = "26.611.82"
$ws.Range("E18").Value = "  -6.92%  "
$ws.Range("D19").Value = "'4.996"
$ws.Range("E19").Value = "  -7.89%  "
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").Value = "'10.94"
$ws.Range("E21").Value = "  -4.48%  "
$ws.Range("D22").Value = "'186.72"
$ws.Range("E22").Value = "  -10.60%  "
$ws.Range("D23").Value = "'6.264"
$ws.Range("E23").Value = "  -8.12%  "
$ws.Range("D24").Value = "'1.007"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").Value = "'144.88"
$ws.Range("E25").Value = "  -5.17%  "
$ws.Range("D26").Value = "'7.492"
$ws.Range("E26").Value = "  -7.94%  "
$ws.Range("D27").Value = "'0.1172"
$ws.Range("E27").Value = "  -6.93%  "
$ws.Range("D28").Value = "'15.81"
$ws.Range("E28").Value = "  -3.47%  "
$ws.Range("D29").Value = "'1.340"
$ws.Range("E29").Value = "  -5.45%  "
$ws.Range("D30").Value = "'0.05735"
$ws.Range("E30").Value = "  -7.38%  "
$ws.Range("D31").Value = "'1.337"
$ws.Range("E31").Value = "  -5.89%  "
$ws.Range("D32").Value = "'3.521"
$ws.Range("E32").Value = "  -6.89%  "
$ws.Range("D33").Value = "'3.506"
$ws.Range("E33").Value = "  -7.85%  "
$ws.Range("D34").Value = "'1.651"
$ws.Range("E34").Value = "  -4.91%  "
$ws.Range("D35").Value = "'1.024"
$ws.Range("E35").Value = "  -1.96%  "
$ws.Range("D36").Value = "'0.5997"
$ws.Range("E36").Value = "  -6.12%  "
$ws.Range("D37").Value = "'2.363"
$ws.Range("E37").Value = "  -5.54%  "
$ws.Range("D38").Value = "'2.684"
$ws.Range("E38").Value = "  -1.01%  "
$ws.Range("D39").Value = "'0.01623"
$ws.Range("E39").Value = "  -4.16%  "
$ws.Range("D40").Value = "1.091.24"
$ws.Range("E40").Value = "  -4.31%  "
$ws.Range("D41").Value = "'0.8630"
$ws.Range("E41").Value = "  -1.55%  "
$ws.Range("D42").Value = "'5.831"
$ws.Range("E42").Value = "  -9.02%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").Value = "'99.80"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").Value = "1.851.54"
$ws.Range("E45").Value = "  -4.82%  "
$ws.Range("E46").Value = "  +6.81%  "
$ws.Range("D47").Value = "'56.45"
$ws.Range("E47").Value = "  -5.73%  "
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("D49").Value = "'8.118"
$ws.Range("E49").Value = "  -2.85%  "
$ws.Range("D50").Value = "'0.05237"
$ws.Range("E50").Value = "  -4.29%  "
$ws.Range("D51").Value = "'0.4320"
$ws.Range("E51").Value = "  -3.51%  "
